$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "corAnsA" answer key in column F (rows 2-9) had its "k"/"l" values
# swapped with the "corAnsF" key in column E. Correct it by toggling every
# F-column value between "k" and "l" for the data rows.
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $current = $cell.Value2
    if ($current -eq "k") {
        $cell.Value = "l"
    } elseif ($current -eq "l") {
        $cell.Value = "k"
    }
}

# Restore the active cell selection left by the author.
$ws.Range("I9").Select()
